$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Update existing D-column figures (re-export of the source data) ---
$dUpdates = @{
    2  = 11767
    3  = 11512
    7  = 11871
    8  = 11159
    12 = 11979
    13 = 10657
    17 = 12023
    18 = 10106
    22 = 12080
    23 = 9676
    27 = 12122
    28 = 9215
    32 = 12162
    33 = 8807
    37 = 12200
    38 = 8283
    42 = 12235
    43 = 7715
    47 = 12260
    48 = 7167
    52 = 12288
    53 = 6525
    57 = 12311
    58 = 5915
    62 = 12322
    63 = 5749
    67 = 12344
    68 = 5565
    72 = 12360
    73 = 5129
    74 = 7231
    75 = 1416
    76 = 5815
    77 = 12383
    78 = 4652
    79 = 7731
    80 = 1472
    81 = 6259
}

foreach ($row in $dUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $dUpdates[$row]
}

# --- Append the new 202504 week (rows 82-86) ---
$newRows = @(
    @{ Row = 82; Variable = "farms_total_count";             Number = 12401 },
    @{ Row = 83; Variable = "farms_to_examine_count";         Number = 4134  },
    @{ Row = 84; Variable = "farms_examined_count";            Number = 8267  },
    @{ Row = 85; Variable = "farms_examined_positive_count";   Number = 1508  },
    @{ Row = 86; Variable = "farms_examined_negative_count";   Number = 6759  }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Copy the format of the row directly above so the new row (esp. the
    # LastDayOfWeek date cell) picks up the same cell style, then overwrite
    # the values.
    $ws.Range("A" + ($r - 1) + ":D" + ($r - 1)).Copy($ws.Range("A" + $r + ":D" + $r))

    $ws.Cells.Item($r, 1).Value = 202504
    $ws.Cells.Item($r, 2).Value = 45683
    $ws.Cells.Item($r, 3).Value = $entry.Variable
    $ws.Cells.Item($r, 4).Value = $entry.Number
}

# --- Restore the view state shown in the final workbook ---
$null = $ws.Range("F18").Select()
